# Generate Report for Handoff
#
# The localization round for the 50ba612b.../c5e9cb3d... files has moved
# on from "handed back" to "ready for handoff" again: the handback status
# text, the handoff priority and a couple of timestamps are refreshed, and
# the zh-cn / de-de tables now flag that the handback file version is
# stale via a new Error Detail message.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$priority = "mt"
$hoXliffDate = "2016-10-20 09:52:13"
$zhHandoffDate = "2016-10-20 09:52:01"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/865e966a73debd1ba21973ece1e28af3a0d96289/e2e/50ba612b-7526-4296-a382-71bb777d8ff3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7a4e2cfaaca073b148b09b1720fa98edb1bdb9b/e2e/50ba612b-7526-4296-a382-71bb777d8ff3.md."

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("G2").Value = $hoXliffDate
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $hoXliffDate

# Columns E/F shrink now that the status text is shorter (autofit-style
# resize performed by the tool that produced this report). The COM layer
# quantizes ColumnWidth to 1/6-character steps, so feed it the input that
# rounds to the closest achievable width to the recorded target
# (17.2159881591797).
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $status
$wsZhCn.Range("E2").Value = $priority
$wsZhCn.Range("H2").Value = $zhHandoffDate
$wsZhCn.Range("P2").Value = $errorDetail
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("E3").Value = $priority
$wsZhCn.Range("H3").Value = $zhHandoffDate

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $status
$wsDeDe.Range("E2").Value = $priority
$wsDeDe.Range("H2").Value = $hoXliffDate
$wsDeDe.Range("P2").Value = $errorDetail
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("E3").Value = $priority
$wsDeDe.Range("H3").Value = $hoXliffDate

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
